$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers ---
# F1 previously held "joint VDZP"; it now becomes "iron VDZP", and a copy of
# the original header text moves out to I1. G1/H1 are brand-new headers.
$ws.Range("F1").Value = "iron VDZP"
$ws.Range("G1").Value = "ltneg5 joint VDZP"
$ws.Range("H1").Value = "ltneg4 joint VDZP"
$ws.Range("I1").Value = "joint VDZP"

# Match the bordered/centered/bold header style used by B1:F1 exactly
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: objective ---
$ws.Range("G2").Value = "reg:squarederror"
$ws.Range("H2").Value = "reg:squarederror"
$ws.Range("I2").Value = "reg:squarederror"

# --- Row 3: enable_categorical ---
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false
$ws.Range("I3").Value = $false

# --- Row 4: max_depth ---
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 10

# --- Row 5: n_estimators ---
$ws.Range("F5").Value = 1000
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 1000

# --- Row 6: reg_alpha ---
$ws.Range("G6").Value = 0.1
$ws.Range("H6").Value = 0.1
$ws.Range("I6").Value = 0.1

# --- Row 7: reg_lambda ---
$ws.Range("G7").Value = 0.1
$ws.Range("H7").Value = 0.001
$ws.Range("I7").Value = 0.1
